$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.795.20'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.30%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.850.65'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.71%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.21'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.00%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.37'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.60%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.854.32'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.54%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.49%  '

$ws.Range('E10').Value = '  -0.80%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.35'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.01%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.456'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.67%  '

$ws.Range('E13').Value = '  +0.37%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.84'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.36%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.495.98'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.59%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.859.70'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.19%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.868.10'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.33%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.07'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +6.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.34'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.00%  '

$ws.Range('E20').Value = '  -1.16%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.94'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.97%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '462.74'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.36%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.728'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.38%  '

$ws.Range('E24').Value = '  -3.70%  '

$ws.Range('E25').Value = '  -1.61%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.24'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.32%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.11'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.87%  '

$ws.Range('E28').Value = '  -0.04%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.98'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.99%  '

$ws.Range('E30').Value = '  +0.17%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.997.39'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.64%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.74'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.06%  '

$ws.Range('E33').Value = '  -2.21%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.02'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.11%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.30'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.91%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.823.91'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.92%  '

$ws.Range('E37').Value = '  -1.89%  '

$ws.Range('E38').Value = '  -1.75%  '

$ws.Range('E39').Value = '  +0.01%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.89'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.24'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.19%  '

$ws.Range('E42').Value = '  +0.20%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.311'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.91%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '425.80'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.29%  '

$ws.Range('E45').Value = '  -0.55%  '

$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '47.15'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.82%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.49'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.76%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000273'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.14%  '

$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '40.55'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.44%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.68'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.26%  '
